$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (originally row 26) and the "SC 92" row
# (originally row 28). Deleting row 26 first shifts "SC 92" up to row 27,
# so delete that row next.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()
